$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 373.2
$ws.Range("I15").Value = 373.2
$ws.Range("K15").Value = 1119.6
$ws.Range("M15").Value = -950.5999999999999
$ws.Range("H17").Value = 1545617.8
$ws.Range("J17").Value = 1685901.6
$ws.Range("L17").Value = 5057704.800000001
$ws.Range("N17").Value = -5058040.800000001
$ws.Range("H61").Value = 118.333336
$ws.Range("I61").Value = 118.333336
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 355.000008
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -183.000008
$ws.Range("H97").Value = 23336.4
$ws.Range("J97").Value = 32059.143
$ws.Range("L97").Value = 96177.429
$ws.Range("N97").Value = -97169.429
$ws.Range("H99").Value = 358.27274
$ws.Range("I99").Value = 345.25
$ws.Range("K99").Value = 1035.75
$ws.Range("M99").Value = 462.25
$ws.Range("H100").Value = 1237.6086
$ws.Range("I100").Value = 1466.4445
$ws.Range("K100").Value = 1466.4445
$ws.Range("M100").Value = -925.4445000000001
$ws.Range("H107").Value = 1043.5333
$ws.Range("I107").Value = 1143.4445
$ws.Range("J107").Value = 893.6667
$ws.Range("K107").Value = 1143.4445
$ws.Range("L107").Value = 893.6667
$ws.Range("M107").Value = 776.5554999999999
$ws.Range("N107").Value = -4733.6667
$ws.Range("H138").Value = 6669284.5
$ws.Range("I138").Value = 1278.6786
$ws.Range("J138").Value = 10641713
$ws.Range("K138").Value = 3836.0358
$ws.Range("L138").Value = 31925139
$ws.Range("M138").Value = 1303.9642
$ws.Range("N138").Value = -31935419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 669.4
$ws.Range("I4").Value = 211.75
$ws.Range("K4").Value = 211.75
$ws.Range("M4").Value = -95.75
$ws.Range("H32").Value = 10216.743
$ws.Range("I32").Value = 6638.2856
$ws.Range("K32").Value = 6638.2856
$ws.Range("M32").Value = -6351.2856
$ws.Range("H45").Value = 11109.9
$ws.Range("I45").Value = 12122.111
$ws.Range("K45").Value = 12122.111
$ws.Range("M45").Value = -11745.111
$ws.Range("H61").Value = 5404.1665
$ws.Range("I61").Value = 4812.5
$ws.Range("K61").Value = 4812.5
$ws.Range("M61").Value = -4600.5
$ws.Range("H102").Value = 2830.9
$ws.Range("I102").Value = 2330.2856
$ws.Range("J102").Value = 3999
$ws.Range("K102").Value = 2330.2856
$ws.Range("L102").Value = 3999
$ws.Range("M102").Value = -708.2856000000002
$ws.Range("N102").Value = -7243
$ws.Range("H110").Value = 6648.875
$ws.Range("I110").Value = 7374.8823
$ws.Range("J110").Value = 4885.7144
$ws.Range("K110").Value = 7374.8823
$ws.Range("L110").Value = 4885.7144
$ws.Range("M110").Value = -5329.8823
$ws.Range("N110").Value = -8975.714400000001
$ws.Range("H122").Value = 1498.1364
$ws.Range("I122").Value = 1225
$ws.Range("K122").Value = 3675
$ws.Range("M122").Value = -1225
$ws.Range("H136").Value = 5404.1665
$ws.Range("I136").Value = 4812.5
$ws.Range("K136").Value = 14437.5
$ws.Range("M136").Value = -11887.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4691.769
$ws.Range("I20").Value = 3711.3333
$ws.Range("J20").Value = 5532.143
$ws.Range("K20").Value = 3711.3333
$ws.Range("L20").Value = 5532.143
$ws.Range("M20").Value = -3464.3333
$ws.Range("N20").Value = -6026.143
$ws.Range("H96").Value = 14428
$ws.Range("I96").Value = 14428
$ws.Range("K96").Value = 14428
$ws.Range("M96").Value = -11682
$ws.Range("H105").Value = 1930.4445
$ws.Range("I105").Value = 1996.75
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 1996.75
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -249.75
$ws.Range("N105").Value = -4894
$ws.Range("H107").Value = 931.25
$ws.Range("I107").Value = 852.88464
$ws.Range("K107").Value = 852.88464
$ws.Range("M107").Value = 1067.11536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2790.4211
$ws.Range("I16").Value = 2645.7693
$ws.Range("K16").Value = 2645.7693
$ws.Range("M16").Value = -2358.7693
$ws.Range("H60").Value = 16098.75
$ws.Range("J60").Value = 17951.5
$ws.Range("L60").Value = 17951.5
$ws.Range("N60").Value = -18973.5
$ws.Range("H99").Value = 3127.75
$ws.Range("I99").Value = 2756
$ws.Range("J99").Value = 3499.5
$ws.Range("K99").Value = 2756
$ws.Range("L99").Value = 3499.5
$ws.Range("M99").Value = -1258
$ws.Range("N99").Value = -6495.5
$ws.Range("H113").Value = 2790.4211
$ws.Range("I113").Value = 2645.7693
$ws.Range("K113").Value = 2645.7693
$ws.Range("M113").Value = -475.7692999999999
$ws.Range("H122").Value = 1749.8125
$ws.Range("I122").Value = 1360.8889
$ws.Range("K122").Value = 4082.6667
$ws.Range("M122").Value = -1632.6667
$ws.Range("H126").Value = 3127.75
$ws.Range("I126").Value = 2756
$ws.Range("J126").Value = 3499.5
$ws.Range("K126").Value = 8268
$ws.Range("L126").Value = 10498.5
$ws.Range("M126").Value = -5798
$ws.Range("N126").Value = -15438.5
$ws.Range("H132").Value = 3170.6667
$ws.Range("I132").Value = 2972.9285
$ws.Range("K132").Value = 8918.7855
$ws.Range("M132").Value = -6388.7855
$ws.Range("H133").Value = 57113
$ws.Range("J133").Value = 57113
$ws.Range("L133").Value = 57113
$ws.Range("N133").Value = -62173
$ws.Range("H134").Value = 16351.156
$ws.Range("I134").Value = 9083.666999999999
$ws.Range("J134").Value = 55595.6
$ws.Range("K134").Value = 27251.001
$ws.Range("L134").Value = 166786.8
$ws.Range("M134").Value = -24716.001
$ws.Range("N134").Value = -171856.8
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0
$ws.Range("H141").Value = 423441.66
$ws.Range("I141").Value = 69999.5
$ws.Range("J141").Value = 600162.75
$ws.Range("K141").Value = 69999.5
$ws.Range("L141").Value = 600162.75
$ws.Range("M141").Value = -64819.5
$ws.Range("N141").Value = -610522.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3075.75
$ws.Range("I34").Value = 150
$ws.Range("J34").Value = 6001.5
$ws.Range("K34").Value = 450
$ws.Range("L34").Value = 18004.5
$ws.Range("M34").Value = -366
$ws.Range("N34").Value = -18172.5
$ws.Range("H55").Value = 1057.6
$ws.Range("J55").Value = 2200
$ws.Range("L55").Value = 6600
$ws.Range("N55").Value = -6954

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1180.7
$ws.Range("I97").Value = 1254.625
$ws.Range("J97").Value = 885
$ws.Range("K97").Value = 1254.625
$ws.Range("L97").Value = 885
$ws.Range("M97").Value = -758.625
$ws.Range("N97").Value = -1877
$ws.Range("H102").Value = 76923710
$ws.Range("I102").Value = 683.8333
$ws.Range("J102").Value = 1000000000
$ws.Range("K102").Value = 683.8333
$ws.Range("L102").Value = 1000000000
$ws.Range("M102").Value = 938.1667
$ws.Range("N102").Value = -1000003244
$ws.Range("H132").Value = 4182.905
$ws.Range("I132").Value = 3437.3076
$ws.Range("J132").Value = 5394.5
$ws.Range("K132").Value = 10311.9228
$ws.Range("L132").Value = 16183.5
$ws.Range("M132").Value = -7781.9228
$ws.Range("N132").Value = -21243.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1648.826
$ws.Range("I22").Value = 773
$ws.Range("J22").Value = 2451.6667
$ws.Range("K22").Value = 773
$ws.Range("L22").Value = 2451.6667
$ws.Range("M22").Value = -478
$ws.Range("N22").Value = -3041.6667
$ws.Range("H27").Value = 1648.826
$ws.Range("I27").Value = 773
$ws.Range("J27").Value = 2451.6667
$ws.Range("K27").Value = 773
$ws.Range("L27").Value = 2451.6667
$ws.Range("M27").Value = -666
$ws.Range("N27").Value = -2665.6667
$ws.Range("H40").Value = 3515.8572
$ws.Range("I40").Value = 2685.1667
$ws.Range("J40").Value = 8500
$ws.Range("K40").Value = 2685.1667
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = -2549.1667
$ws.Range("N40").Value = -8772
$ws.Range("H88").Value = 12500
$ws.Range("I88").Value = 12500
$ws.Range("K88").Value = 12500
$ws.Range("M88").Value = -12072
$ws.Range("H91").Value = 12500
$ws.Range("I91").Value = 12500
$ws.Range("K91").Value = 12500
$ws.Range("M91").Value = -11018
$ws.Range("H132").Value = 6017.3335
$ws.Range("I132").Value = 5750
$ws.Range("K132").Value = 17250
$ws.Range("M132").Value = -14720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 45966.332
$ws.Range("I69").Value = 44900
$ws.Range("K69").Value = 44900
$ws.Range("M69").Value = -44151
$ws.Range("H72").Value = 45966.332
$ws.Range("I72").Value = 44900
$ws.Range("K72").Value = 134700
$ws.Range("M72").Value = -130956
$ws.Range("H100").Value = 590.8461
$ws.Range("I100").Value = 585.2632
$ws.Range("K100").Value = 1170.5264
$ws.Range("M100").Value = -629.5264
$ws.Range("H122").Value = 2904.2222
$ws.Range("I122").Value = 2858.1667
$ws.Range("K122").Value = 8574.500100000001
$ws.Range("M122").Value = -6124.500100000001
$ws.Range("H126").Value = 1861.5652
$ws.Range("J126").Value = 1947.5714
$ws.Range("L126").Value = 5842.7142
$ws.Range("N126").Value = -10782.7142
